{"js": "const replacements = [\n  [\"2025-06-02 Monday\", \"2025-06-03 Tuesday\"],\n  [\"733\u00f77=\", \"513\u00f78=\"],\n  [\"742\u00f76=\", \"810\u00f74=\"],\n  [\"342\u00f75=\", \"876\u00f74=\"],\n  [\"622\u00f75=\", \"307\u00f79=\"],\n  [\"557\u00f77=\", \"238\u00f72=\"],\n  [\"400\u00f74=\", \"531\u00f79=\"],\n  [\"527\u00f77=\", \"734\u00f76=\"],\n  [\"456\u00f78=\", \"184\u00f73=\"],\n  [\"328\u00f77=\", \"981\u00f73=\"],\n  [\"233\u00f76=\", \"898\u00f76=\"],\n  [\"188\u00f73=\", \"170\u00f75=\"],\n  [\"925\u00f72=\", \"711\u00f73=\"],\n  [\"296\u00f79=\", \"550\u00f76=\"],\n  [\"465\u00f77=\", \"143\u00f79=\"],\n  [\"747\u00f78=\", \"924\u00f78=\"],\n  [\"705\u00f78=\", \"327\u00f79=\"],\n  [\"491\u00f79=\", \"247\u00f72=\"],\n  [\"920\u00f75=\", \"840\u00f79=\"],\n  [\"122\u00f79=\", \"513\u00f76=\"],\n  [\"323\u00f73=\", \"805\u00f79=\"],\n  [\"707\u00f79=\", \"200\u00f77=\"],\n  [\"268\u00f72=\", \"155\u00f77=\"],\n  [\"733\u00f72=\", \"400\u00f72=\"],\n  [\"147\u00f79=\", \"191\u00f72=\"],\n  [\"342\u00f73=\", \"627\u00f78=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-06-02 Monday\", \"2025-06-03 Tuesday\"),\n    @(\"733\u00f77=\", \"513\u00f78=\"),\n    @(\"742\u00f76=\", \"810\u00f74=\"),\n    @(\"342\u00f75=\", \"876\u00f74=\"),\n    @(\"622\u00f75=\", \"307\u00f79=\"),\n    @(\"557\u00f77=\", \"238\u00f72=\"),\n    @(\"400\u00f74=\", \"531\u00f79=\"),\n    @(\"527\u00f77=\", \"734\u00f76=\"),\n    @(\"456\u00f78=\", \"184\u00f73=\"),\n    @(\"328\u00f77=\", \"981\u00f73=\"),\n    @(\"233\u00f76=\", \"898\u00f76=\"),\n    @(\"188\u00f73=\", \"170\u00f75=\"),\n    @(\"925\u00f72=\", \"711\u00f73=\"),\n    @(\"296\u00f79=\", \"550\u00f76=\"),\n    @(\"465\u00f77=\", \"143\u00f79=\"),\n    @(\"747\u00f78=\", \"924\u00f78=\"),\n    @(\"705\u00f78=\", \"327\u00f79=\"),\n    @(\"491\u00f79=\", \"247\u00f72=\"),\n    @(\"920\u00f75=\", \"840\u00f79=\"),\n    @(\"122\u00f79=\", \"513\u00f76=\"),\n    @(\"323\u00f73=\", \"805\u00f79=\"),\n    @(\"707\u00f79=\", \"200\u00f77=\"),\n    @(\"268\u00f72=\", \"155\u00f77=\"),\n    @(\"733\u00f72=\", \"400\u00f72=\"),\n    @(\"147\u00f79=\", \"191\u00f72=\"),\n    @(\"342\u00f73=\", \"627\u00f78=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)\n}\n"}
